$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), row 2
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = ""

# Sheet ALC (index 1), row 41
$ws = $wb.Worksheets.Item(1)
$ws.Range("H41").Value = 323.82352
$ws.Range("J41").Value = 337.8125
$ws.Range("L41").Value = 337.8125
$ws.Range("N41").Value = -1217.8125

# Sheet ALC (index 1), row 64
$ws = $wb.Worksheets.Item(1)
$ws.Range("H64").Value = 4166.6665
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 4533.3335
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 4533.3335
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -5029.3335

# Sheet ALC (index 1), row 67
$ws = $wb.Worksheets.Item(1)
$ws.Range("H67").Value = 4166.6665
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 4533.3335
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 4533.3335
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -6249.3335

# Sheet ALC (index 1), row 112
$ws = $wb.Worksheets.Item(1)
$ws.Range("H112").Value = 2600.139
$ws.Range("J112").Value = 2617.2856
$ws.Range("L112").Value = 7851.8568
$ws.Range("N112").Value = -10067.8568

# Sheet ALC (index 1), row 124
$ws = $wb.Worksheets.Item(1)
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""

# Sheet ALC (index 1), row 125
$ws = $wb.Worksheets.Item(1)
$ws.Range("H125").Value = 530.625
$ws.Range("I125").Value = 198.75
$ws.Range("J125").Value = 862.5
$ws.Range("K125").Value = 1788.75
$ws.Range("L125").Value = 7762.5
$ws.Range("M125").Value = 671.25
$ws.Range("N125").Value = -12682.5

# Sheet ALC (index 1), row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 2671.3447
$ws.Range("I132").Value = 2671.3447
$ws.Range("K132").Value = 8014.034100000001
$ws.Range("M132").Value = -5484.034100000001

# Sheet ALC (index 1), row 135
$ws = $wb.Worksheets.Item(1)
$ws.Range("H135").Value = 50014590
$ws.Range("I135").Value = 991.4286
$ws.Range("J135").Value = 166713000
$ws.Range("K135").Value = 8922.857399999999
$ws.Range("L135").Value = 1500417000
$ws.Range("M135").Value = -6387.857399999999
$ws.Range("N135").Value = -1500422070

# Sheet ALC (index 1), row 137
$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value = 314911.62
$ws.Range("I137").Value = 502511.88
$ws.Range("K137").Value = 1507535.64
$ws.Range("M137").Value = -1504985.64

# Sheet ALC (index 1), row 141
$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value = 3512.8572
$ws.Range("I141").Value = 3522.5
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 10567.5
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = -5387.5
$ws.Range("N141").Value = -20860

# Sheet ARM (index 2), row 39
$ws = $wb.Worksheets.Item(2)
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""

# Sheet ARM (index 2), row 45
$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 2363.9363
$ws.Range("I45").Value = 1605.5
$ws.Range("J45").Value = 3302.9524
$ws.Range("K45").Value = 1605.5
$ws.Range("L45").Value = 3302.9524
$ws.Range("M45").Value = -1228.5
$ws.Range("N45").Value = -4056.9524

# Sheet ARM (index 2), row 74
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 2213.611
$ws.Range("I74").Value = 1848.5385
$ws.Range("J74").Value = 3162.8
$ws.Range("K74").Value = 1848.5385
$ws.Range("L74").Value = 3162.8
$ws.Range("M74").Value = -974.5385000000001
$ws.Range("N74").Value = -4910.8

# Sheet ARM (index 2), row 77
$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 2213.611
$ws.Range("I77").Value = 1848.5385
$ws.Range("J77").Value = 3162.8
$ws.Range("K77").Value = 9242.692500000001
$ws.Range("L77").Value = 15814
$ws.Range("M77").Value = -4874.692500000001
$ws.Range("N77").Value = -24550

# Sheet ARM (index 2), row 128
$ws = $wb.Worksheets.Item(2)
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

# Sheet ARM (index 2), row 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 22568.48
$ws.Range("I132").Value = 1919.7142
$ws.Range("J132").Value = 48848.727
$ws.Range("K132").Value = 5759.142599999999
$ws.Range("L132").Value = 146546.181
$ws.Range("M132").Value = -3229.142599999999
$ws.Range("N132").Value = -151606.181

# Sheet BSM (index 3), row 20
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1602.8889
$ws.Range("I20").Value = 1332.2858
$ws.Range("J20").Value = 2550
$ws.Range("K20").Value = 1332.2858
$ws.Range("L20").Value = 2550
$ws.Range("M20").Value = -1085.2858
$ws.Range("N20").Value = -3044

# Sheet BSM (index 3), row 107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 1184.7646
$ws.Range("I107").Value = 767.2727
$ws.Range("J107").Value = 1950.1666
$ws.Range("K107").Value = 767.2727
$ws.Range("L107").Value = 1950.1666
$ws.Range("M107").Value = 1152.7273
$ws.Range("N107").Value = -5790.1666

# Sheet BSM (index 3), row 137
$ws = $wb.Worksheets.Item(3)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# Sheet BSM (index 3), row 138
$ws = $wb.Worksheets.Item(3)
$ws.Range("H138").Value = 49980
$ws.Range("J138").Value = 49980
$ws.Range("L138").Value = 49980
$ws.Range("N138").Value = -60260

# Sheet CRP (index 4), row 7
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 116
$ws.Range("I7").Value = 100.333336
$ws.Range("K7").Value = 100.333336
$ws.Range("M7").Value = 12.666664

# Sheet CRP (index 4), row 22
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 343.33334
$ws.Range("I22").Value = 372
$ws.Range("K22").Value = 372
$ws.Range("M22").Value = -22

# Sheet CRP (index 4), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 14723.303
$ws.Range("J31").Value = 4915.5
$ws.Range("L31").Value = 4915.5
$ws.Range("N31").Value = -5505.5

# Sheet CRP (index 4), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 14723.303
$ws.Range("J34").Value = 4915.5
$ws.Range("L34").Value = 4915.5
$ws.Range("N34").Value = -5319.5

# Sheet CRP (index 4), row 105
$ws = $wb.Worksheets.Item(4)
$ws.Range("H105").Value = 15625829
$ws.Range("J105").Value = 1670.3334
$ws.Range("L105").Value = 1670.3334
$ws.Range("N105").Value = -5164.3334

# Sheet CRP (index 4), row 132
$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 30555.945
$ws.Range("I132").Value = 34466.332
$ws.Range("J132").Value = 11004
$ws.Range("K132").Value = 103398.996
$ws.Range("L132").Value = 33012
$ws.Range("M132").Value = -100868.996
$ws.Range("N132").Value = -38072

# Sheet CUL (index 5), row 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 1118.1875
$ws.Range("I5").Value = 929.6923
$ws.Range("K5").Value = 2789.0769
$ws.Range("M5").Value = -2677.0769

# Sheet CUL (index 5), row 37
$ws = $wb.Worksheets.Item(5)
$ws.Range("H37").Value = 22792636
$ws.Range("J37").Value = 22792636
$ws.Range("L37").Value = 68377908
$ws.Range("N37").Value = -68378132

# Sheet CUL (index 5), row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 753.67
$ws.Range("J131").Value = 753.67
$ws.Range("L131").Value = 2261.01
$ws.Range("N131").Value = -12341.01

# Sheet CUL (index 5), row 135
$ws = $wb.Worksheets.Item(5)
$ws.Range("H135").Value = 1118.1875
$ws.Range("I135").Value = 929.6923
$ws.Range("K135").Value = 8367.2307
$ws.Range("M135").Value = -5832.2307

# Sheet CUL (index 5), row 136
$ws = $wb.Worksheets.Item(5)
$ws.Range("H136").Value = 1418.2222
$ws.Range("I136").Value = 966.375
$ws.Range("J136").Value = 5033
$ws.Range("K136").Value = 2899.125
$ws.Range("L136").Value = 15099
$ws.Range("M136").Value = 2200.875
$ws.Range("N136").Value = -25299

# Sheet CUL (index 5), row 137
$ws = $wb.Worksheets.Item(5)
$ws.Range("H137").Value = 11906804
$ws.Range("I137").Value = 987.8333
$ws.Range("J137").Value = 15153845
$ws.Range("K137").Value = 2963.4999
$ws.Range("L137").Value = 45461535
$ws.Range("M137").Value = 2136.5001
$ws.Range("N137").Value = -45471735

# Sheet CUL (index 5), row 138
$ws = $wb.Worksheets.Item(5)
$ws.Range("H138").Value = 1911.625
$ws.Range("I138").Value = 1812
$ws.Range("K138").Value = 5436
$ws.Range("M138").Value = -296

# Sheet CUL (index 5), row 139
$ws = $wb.Worksheets.Item(5)
$ws.Range("H139").Value = 1649.8182
$ws.Range("I139").Value = 1649.8182
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 4949.4546
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 190.5454
$ws.Range("N139").Value = ""

# Sheet GSM (index 6), row 126
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 3941.3333
$ws.Range("I126").Value = 3221.875
$ws.Range("J126").Value = 4618.4707
$ws.Range("K126").Value = 9665.625
$ws.Range("L126").Value = 13855.4121
$ws.Range("M126").Value = -7195.625
$ws.Range("N126").Value = -18795.4121

# Sheet GSM (index 6), row 132
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 78961.8
$ws.Range("I132").Value = 75446.28999999999
$ws.Range("J132").Value = 87164.664
$ws.Range("K132").Value = 226338.87
$ws.Range("L132").Value = 261493.992
$ws.Range("M132").Value = -223808.87
$ws.Range("N132").Value = -266553.992

# Sheet LTW (index 7), row 61
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 3928.25
$ws.Range("I61").Value = 2146.2354
$ws.Range("J61").Value = 6682.273
$ws.Range("K61").Value = 2146.2354
$ws.Range("L61").Value = 6682.273
$ws.Range("M61").Value = -1944.2354
$ws.Range("N61").Value = -7086.273

# Sheet LTW (index 7), row 113
$ws = $wb.Worksheets.Item(7)
$ws.Range("H113").Value = 3928.25
$ws.Range("I113").Value = 2146.2354
$ws.Range("J113").Value = 6682.273
$ws.Range("K113").Value = 2146.2354
$ws.Range("L113").Value = 6682.273
$ws.Range("M113").Value = 23.76459999999997
$ws.Range("N113").Value = -11022.273

# Sheet LTW (index 7), row 132
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 636833.75
$ws.Range("I132").Value = 863031.3
$ws.Range("J132").Value = 3480.6
$ws.Range("K132").Value = 2589093.9
$ws.Range("L132").Value = 10441.8
$ws.Range("M132").Value = -2586563.9
$ws.Range("N132").Value = -15501.8

# Sheet LTW (index 7), row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 38539.5
$ws.Range("I136").Value = 47413.91
$ws.Range("K136").Value = 142241.73
$ws.Range("M136").Value = -139691.73

# Sheet WVR (index 8), row 4
$ws = $wb.Worksheets.Item(8)
$ws.Range("H4").Value = 8500
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -1887
$ws.Range("N4").Value = -15226

# Sheet WVR (index 8), row 81
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 2157.1428
$ws.Range("J81").Value = 3200
$ws.Range("L81").Value = 6400
$ws.Range("N81").Value = -8522

# Sheet WVR (index 8), row 84
$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value = 2157.1428
$ws.Range("J84").Value = 3200
$ws.Range("L84").Value = 32000
$ws.Range("N84").Value = -42608

# Sheet WVR (index 8), row 113
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 1931515.5
$ws.Range("I113").Value = 1199.091
$ws.Range("J113").Value = 9009342
$ws.Range("K113").Value = 3597.273
$ws.Range("L113").Value = 27028026
$ws.Range("M113").Value = -1427.273
$ws.Range("N113").Value = -27032366

# Sheet WVR (index 8), row 126
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 777.2381
$ws.Range("I126").Value = 677.75
$ws.Range("J126").Value = 1095.6
$ws.Range("K126").Value = 2033.25
$ws.Range("L126").Value = 3286.8
$ws.Range("M126").Value = 436.75
$ws.Range("N126").Value = -8226.799999999999

# Sheet WVR (index 8), row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 2796.7368
$ws.Range("I132").Value = 2449.9167
$ws.Range("J132").Value = 3391.2856
$ws.Range("K132").Value = 7349.750100000001
$ws.Range("L132").Value = 10173.8568
$ws.Range("M132").Value = -4819.750100000001
$ws.Range("N132").Value = -15233.8568

# Sheet WVR (index 8), row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 949894.75
$ws.Range("I136").Value = 1344875.8
$ws.Range("K136").Value = 4034627.4
$ws.Range("M136").Value = -4032077.4
